# Commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
# This regenerates the pitcher-log worksheet's column G ("K" = strikeouts) values,
# which replace the previous "Strike#"-derived figures with freshly calculated/simulated
# strikeout counts (s_vals). Only column G (the 7th column) changes; every other column
# -- date, TB, PC, dS0, dSF, IP, I0, IF -- is left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values keyed by row number, as produced by the regenerated s_vals calc.
$newK = @{
    2 = 3
    3 = 2
    4 = 3
    5 = 5
    6 = 1
    7 = 2
    8 = 1
    9 = 3
    10 = 2
    11 = 0
    13 = 0
    14 = 0
    15 = 2
    16 = 1
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 2
    23 = 0
    24 = 2
    25 = 1
    26 = 2
    27 = 3
    28 = 3
    29 = 0
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 2
    35 = 0
    36 = 2
    37 = 0
    38 = 0
    39 = 1
    41 = 1
    42 = 1
    43 = 0
    44 = 2
    45 = 0
    46 = 3
    47 = 0
    48 = 1
    49 = 1
    50 = 1
    51 = 1
    52 = 1
    53 = 0
    54 = 0
    55 = 3
    56 = 1
    57 = 2
    58 = 1
    59 = 2
    60 = 0
    61 = 1
    62 = 0
    63 = 2
    64 = 1
    65 = 0
    66 = 1
    67 = 2
    68 = 0
    69 = 3
    70 = 0
    71 = 1
    72 = 1
    73 = 1
    74 = 1
    75 = 1
    76 = 1
    77 = 3
    78 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}

